# Update "想去人数" (want-to-go count) values in column F for sheets
# "展览" and "全部类型". Both sheets hold the same list of exhibitions,
# so the same row/value updates apply to both.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 205
    3  = 248
    4  = 266
    6  = 257
    7  = 6105
    8  = 45
    11 = 61
    15 = 417
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
